$d = $word.ActiveDocument

# Fix the typo "performa" -> "proforma" in "coding details performa.docx"
$d.Content.Find.Execute("coding details performa.docx", $true, $false, $false, $false, $false,
                         $true, 1, $false, "coding details proforma.docx", 2)
